# Re-applies the betexplorer scraper refresh for
# serbia/prva-liga/2023-2024: a handful of fixtures that already existed
# in the sheet got their home/away (and odds) data corrected/reordered,
# and two newly-scraped fixtures (round continuing into November) were
# appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: swap the match-data columns (F..V) between two existing rows.
# Columns A..E (Indice/pais/torneio/temporada/data_partida) are left
# untouched - only the fixture-specific data moves.
# ---------------------------------------------------------------------
function Swap-MatchRows($sheet, $row1, $row2) {
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $addr1 = $col + $row1
        $addr2 = $col + $row2
        $tmp = $sheet.Range($addr1).Value2
        $sheet.Range($addr1).Value2 = $sheet.Range($addr2).Value2
        $sheet.Range($addr2).Value2 = $tmp
    }
}

# Pairs of rows whose fixture data was swapped.
Swap-MatchRows $ws 3 4
Swap-MatchRows $ws 5 6
Swap-MatchRows $ws 66 67
Swap-MatchRows $ws 82 83
Swap-MatchRows $ws 91 92

# ---------------------------------------------------------------------
# Rows 97..101: the five fixtures were re-ordered (a rotation), not a
# plain pairwise swap. Capture all five values per column first, then
# write them back in the new order so the read doesn't clobber data
# still needed for a later write.
# ---------------------------------------------------------------------
$rotCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $rotCols) {
    $v97 = $ws.Range($col + "97").Value2
    $v98 = $ws.Range($col + "98").Value2
    $v99 = $ws.Range($col + "99").Value2
    $v100 = $ws.Range($col + "100").Value2
    $v101 = $ws.Range($col + "101").Value2

    $ws.Range($col + "97").Value2 = $v101
    $ws.Range($col + "98").Value2 = $v100
    $ws.Range($col + "99").Value2 = $v97
    $ws.Range($col + "100").Value2 = $v98
    $ws.Range($col + "101").Value2 = $v99
}

# ---------------------------------------------------------------------
# Two brand-new fixtures appended at the end (rows 106 & 107), growing
# the used range from A1:V105 to A1:V107. Clone the formatting of the
# last existing data row (105) first so the new Indice cell keeps the
# bordered/centered style and the date cell keeps its date/time number
# format, then fill in the actual values.
# ---------------------------------------------------------------------
$ws.Range("A105:V105").Copy() | Out-Null
$ws.Range("A106:V106").PasteSpecial(-4122) | Out-Null
$ws.Range("A107:V107").PasteSpecial(-4122) | Out-Null

$ws.Range("A106").Value2 = 105
$ws.Range("B106").Value2 = "serbia"
$ws.Range("C106").Value2 = "prva-liga"
$ws.Range("D106").Value2 = "2023-2024"
$ws.Range("E106").Value2 = 45235.70833333334
$ws.Range("F106").Value2 = "Dubocica"
$ws.Range("G106").Value2 = 1
$ws.Range("H106").Value2 = "Jedinstvo U."
$ws.Range("I106").Value2 = 0
$ws.Range("J106").Value2 = 2.36
$ws.Range("K106").Value2 = "04/11/2023 05:12"
$ws.Range("L106").Value2 = 2.39
$ws.Range("M106").Value2 = "05/11/2023 16:38"
$ws.Range("N106").Value2 = 2.67
$ws.Range("O106").Value2 = "04/11/2023 05:12"
$ws.Range("P106").Value2 = 2.8
$ws.Range("Q106").Value2 = "05/11/2023 16:38"
$ws.Range("R106").Value2 = 2.88
$ws.Range("S106").Value2 = "04/11/2023 05:12"
$ws.Range("T106").Value2 = 3.04
$ws.Range("U106").Value2 = "05/11/2023 16:38"
$ws.Range("V106").Value2 = "https://www.betexplorer.com/football/serbia/prva-liga/dubocica-jedinstvo-ub/MBO1j6bH/"

$ws.Range("A107").Value2 = 106
$ws.Range("B107").Value2 = "serbia"
$ws.Range("C107").Value2 = "prva-liga"
$ws.Range("D107").Value2 = "2023-2024"
$ws.Range("E107").Value2 = 45235.70833333334
$ws.Range("F107").Value2 = "Kolubara"
$ws.Range("G107").Value2 = 0
$ws.Range("H107").Value2 = "Mladost GAT"
$ws.Range("I107").Value2 = 2
$ws.Range("J107").Value2 = 1.99
$ws.Range("K107").Value2 = "04/11/2023 05:13"
$ws.Range("L107").Value2 = 2.15
$ws.Range("M107").Value2 = "05/11/2023 16:44"
$ws.Range("N107").Value2 = 2.81
$ws.Range("O107").Value2 = "04/11/2023 05:13"
$ws.Range("P107").Value2 = 2.78
$ws.Range("Q107").Value2 = "05/11/2023 16:44"
$ws.Range("R107").Value2 = 3.46
$ws.Range("S107").Value2 = "04/11/2023 05:13"
$ws.Range("T107").Value2 = 3.57
$ws.Range("U107").Value2 = "05/11/2023 16:44"
$ws.Range("V107").Value2 = "https://www.betexplorer.com/football/serbia/prva-liga/kolubara-mladost-gat/UuORSIHh/"
